$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7-11 (no longer needed)
$ws.Range("A7:B11").EntireRow.Delete()

# Update values for rows 2-6 (columns A and B)
$values = @(
    @(0, 82),
    @(2, 78),
    @(4, 74),
    @(1, 38),
    @(3, 26)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}
